# Update cryptocurrency price (D) and 1h volume change (E) columns
# to reflect the latest scraped values from coinranking.com.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$c = $ws.Range("D2")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '30.288.16'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style
$ws.Range("E2").Value = '  +0.03%  '

# Row 3
$c = $ws.Range("D3")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '1.867.38'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style
$ws.Range("E3").Value = '  +0.18%  '

# Row 4
$c = $ws.Range("D4")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '1.000'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style
$ws.Range("E4").Value = '  -0.07%  '

# Row 5
$c = $ws.Range("D5")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '234.79'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style
$ws.Range("E5").Value = '  -0.64%  '

# Row 6
$ws.Range("E6").Value = '  -0.04%  '

# Row 7
$ws.Range("E7").Value = '  -0.17%  '

# Row 8
$c = $ws.Range("D8")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '0.2858'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style
$ws.Range("E8").Value = '  -1.57%  '

# Row 9
$c = $ws.Range("D9")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '0.06573'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style
$ws.Range("E9").Value = '  +0.51%  '

# Row 10
$c = $ws.Range("D10")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '21.37'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style
$ws.Range("E10").Value = '  -2.70%  '

# Row 11
$c = $ws.Range("D11")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '0.07823'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style
$ws.Range("E11").Value = '  -1.49%  '

# Row 12
$c = $ws.Range("D12")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '96.78'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style
$ws.Range("E12").Value = '  -1.21%  '

# Row 13
$c = $ws.Range("D13")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '1.864.13'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style
$ws.Range("E13").Value = '  -0.02%  '

# Row 14
$c = $ws.Range("D14")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '0.6966'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style
$ws.Range("E14").Value = '  +2.31%  '

# Row 15
$c = $ws.Range("D15")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '5.094'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style
$ws.Range("E15").Value = '  -1.06%  '

# Row 16
$c = $ws.Range("D16")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '268.17'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style
$ws.Range("E16").Value = '  +2.05%  '

# Row 17
$c = $ws.Range("D17")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '30.396.00'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style
$ws.Range("E17").Value = '  +0.45%  '

# Row 18
$ws.Range("E18").Value = '  +0.37%  '

# Row 19
$c = $ws.Range("D19")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '0.000007656'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style
$ws.Range("E19").Value = '  +2.40%  '

# Row 20
$ws.Range("E20").Value = '  -0.08%  '

# Row 21
$c = $ws.Range("D21")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '2.132.91'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style
$ws.Range("E21").Value = '  +1.10%  '

# Row 22
$c = $ws.Range("D22")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '1.0000'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style

# Row 23
$c = $ws.Range("D23")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '5.241'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style
$ws.Range("E23").Value = '  -0.60%  '

# Row 24
$c = $ws.Range("D24")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '6.176'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style
$ws.Range("E24").Value = '  -0.05%  '

# Row 25
$c = $ws.Range("D25")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '9.486'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style
$ws.Range("E25").Value = '  +3.14%  '

# Row 26
$c = $ws.Range("D26")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '166.56'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style
$ws.Range("E26").Value = '  -0.48%  '

# Row 27
$c = $ws.Range("D27")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '18.87'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style

# Row 28
$c = $ws.Range("D28")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '1.939'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style
$ws.Range("E28").Value = '  -0.60%  '

# Row 29
$ws.Range("E29").Value = '  -2.16%  '

# Row 30
$c = $ws.Range("D30")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '0.09913'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style
$ws.Range("E30").Value = '  -0.08%  '

# Row 31
$c = $ws.Range("D31")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '4.359'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style
$ws.Range("E31").Value = '  +0.21%  '

# Row 32
$ws.Range("E32").Value = '  -0.87%  '

# Row 33
$c = $ws.Range("D33")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '4.048'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style
$ws.Range("E33").Value = '  +0.68%  '

# Row 34
$c = $ws.Range("D34")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '0.04725'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style
$ws.Range("E34").Value = '  +0.20%  '

# Row 35
$ws.Range("E35").Value = '  +0.06%  '

# Row 36
$ws.Range("E36").Value = '  +0.35%  '

# Row 37
$ws.Range("E37").Value = '  +0.38%  '

# Row 38
$ws.Range("E38").Value = '  -0.37%  '

# Row 39
$c = $ws.Range("D39")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '2.752'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style
$ws.Range("E39").Value = '  +5.02%  '

# Row 40
$c = $ws.Range("D40")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '6.316'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style
$ws.Range("E40").Value = '  -0.08%  '

# Row 41
$c = $ws.Range("D41")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '72.76'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style
$ws.Range("E41").Value = '  -1.52%  '

# Row 42
$ws.Range("E42").Value = '  +0.33%  '

# Row 43
$c = $ws.Range("D43")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '0.4171'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style
$ws.Range("E43").Value = '  +0.27%  '

# Row 44
$ws.Range("E44").Value = '  +0.07%  '

# Row 45
$c = $ws.Range("D45")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '0.8343'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style
$ws.Range("E45").Value = '  -0.93%  '

# Row 46
$c = $ws.Range("D46")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '103.05'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style
$ws.Range("E46").Value = '  -0.19%  '

# Row 47
$c = $ws.Range("D47")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '972.10'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style
$ws.Range("E47").Value = '  +2.32%  '

# Row 48
$c = $ws.Range("D48")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '7.108'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style
$ws.Range("E48").Value = '  -0.75%  '

# Row 49
$c = $ws.Range("D49")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '9.110'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style
$ws.Range("E49").Value = '  -1.08%  '

# Row 50
$c = $ws.Range("D50")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '34.48'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style

# Row 51
$c = $ws.Range("D51")
$c.NumberFormat = "@"   # force text so values like "1.000" are not parsed as numbers
$c.Value = '0.05686'
$c.ClearFormats()       # drop the temporary text format, restoring default cell style
$ws.Range("E51").Value = '  +0.38%  '
